# edit.ps1 - "committing updated files for dbBot, url2 and 3 were not
# needed anymore. Edited the word document generator."
#
#  - Sheet1: URL3 column (I) is no longer needed -> deleted entirely.
#    URL2 header (H1) renamed to "Contributer" and every data row's H
#    value becomes "Bella" (the contributor / dbBot result column).
#  - Sheet1: Institution / City data (Temple University / Phladelphia, PA)
#    now filled in for every data row; Date_Created / Date_Expired
#    timestamps refreshed to the latest bot run.
#  - Sheet1 column G widened a bit for the longer "Yes/No, individual..."
#    text.
#  - View bookkeeping: Sheet1 is now the active tab (selection sitting
#    just below the data, at A12); sheet2 is no longer active (its
#    selection covers the full A1:A17 block it was used to generate the
#    Word document from).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("sheet2")

# --- Sheet1: remove the now-unneeded URL3 column, rename URL2 ----------
$ws1.Columns.Item(9).Delete()
$ws1.Range("H1").Value = "Contributer"

# widen column G slightly (closest the engine's column-width grid allows
# to the authored 34.140625 -- the host snaps ColumnWidth to sixths of a
# character, so feed it the pre-image that lands on the nearest bucket)
$ws1.Columns.Item(7).ColumnWidth = 33.33333333

# --- Sheet1: refreshed Date_Created / Date_Expired + new row data ------
$rows = @(
    @{ Row = 2;  C = 43762.61361264891; D = 44128.61361264891 },
    @{ Row = 3;  C = 43762.61364490058; D = 44128.61364490058 },
    @{ Row = 4;  C = 43762.61366567847; D = 44128.61366567847 },
    @{ Row = 5;  C = 43762.61368709357; D = 44128.61368709357 },
    @{ Row = 6;  C = 43762.61370829852; D = 44128.61370829852 },
    @{ Row = 7;  C = 43762.61373041893; D = 44128.61373041893 },
    @{ Row = 8;  C = 43762.61375293967; D = 44128.61375293967 },
    @{ Row = 9;  C = 43762.6137754529;  D = 44128.6137754529  },
    @{ Row = 10; C = 43762.61379878169; D = 44128.61379878169 },
    @{ Row = 11; C = 43762.61382257242; D = 44128.61382257242 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws1.Cells.Item($rowNum, 3).Value = $r.C
    $ws1.Cells.Item($rowNum, 4).Value = $r.D
    $ws1.Cells.Item($rowNum, 5).Value = "Temple University"
    $ws1.Cells.Item($rowNum, 6).Value = "Phladelphia, PA "
    $ws1.Cells.Item($rowNum, 8).Value = "Bella"
}

# --- view bookkeeping ----------------------------------------------------
# sheet2 selection now spans A1:A17 (its last active selection while
# generating the Word doc)
$ws2.Range("A1:A17").Select()

# Sheet1 becomes the active tab, selection resting at A12 just under the
# table
$ws1.Range("A12").Select()

$wb.Save()
